$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.475.83"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.377.45"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'319.76"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'109.50"
$ws.Range("E6").Value = "  -4.52%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'41.25"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").Value = "'0.0922"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "'8.54"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'0.986"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "2.739.98"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "2.383.41"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "45.421.11"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +14.65%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "'0.0000107"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'267.20"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'7.53"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'11.22"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "'22.57"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "'0.0951"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "'37.34"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "'169.48"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'2.84"
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("E38").Value = "  +11.32%  "
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "'98.75"
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").Value = "'70.49"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "1.876.10"
$ws.Range("E44").Value = "  +14.20%  "
$ws.Range("D45").Value = "'12.98"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "'0.228"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "'5.97"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "'84.30"
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("D50").Value = "'112.33"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").Value = "'9.26"
$ws.Range("E51").Value = "  -1.14%  "
